$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting rows 7:13 down to 8:14
$ws.Rows.Item(7).Insert()

# Fill the new row 7 with the weekly data point
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44875
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7500
$ws.Range("P7").Value = 7250
$ws.Range("Q7").Value = '$/bandeja 2 kilos'
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 3625
$ws.Range("T7").Value = 2
